# [FIX] fix email record template,receive setting;
#
# The NOTIFY_SEND_SETTING sheet had three extra "email record template"
# receive-setting rows (issueCreate / issueAssignee / issueSolve) whose
# CODE/NAME/DESCRIPTION/FD_LEVEL columns (E:J) were populated by mistake.
# Clear those columns on rows 17-19, leaving only the RETRY_COUNT /
# IS_SEND_INSTANTLY / IS_MANUAL_RETRY flag columns (K:M) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOTIFY_SEND_SETTING")

$ws.Range("E17:J19").ClearContents()

# Restore the sheet view: scrolled back to the top-left and the formerly
# data-filled block (now emptied) selected, matching the saved workbook view.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("D17:J19").Select()
